# Update the yearly database worksheet:
#  - shift the "twelve months ended" period labels forward by one year
#    (drop 1396/12, add 1401/12 as the newest period)
#  - refresh every data row with the newest figures. Most values simply
#    move one column to the left and a freshly computed figure is added
#    for the new 1401/12 column; the 1400/12 figure for
#    "هزینه حقوق و دستمزد" (row 19) was also recomputed because of the
#    updated read_price algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Period header labels (used in both row 8 and row 24) ----
# Replacing from the newest label down to the oldest lets each
# replacement free up the shared string that the next step needs,
# so the whole workbook ends up referencing a clean, de-duplicated
# set of period labels instead of leaving orphaned/duplicated text.
$ws.Cells.Replace("دوازده ماهه منتهی به 1400/12", "دوازده ماهه منتهی به 1401/12")
$ws.Cells.Replace("دوازده ماهه منتهی به 1399/12", "دوازده ماهه منتهی به 1400/12")
$ws.Cells.Replace("دوازده ماهه منتهی به 1398/12", "دوازده ماهه منتهی به 1399/12")
$ws.Cells.Replace("دوازده ماهه منتهی به 1397/12", "دوازده ماهه منتهی به 1398/12")
$ws.Cells.Replace("دوازده ماهه منتهی به 1396/12", "دوازده ماهه منتهی به 1397/12")

# ---- Data rows ----
# Map of row number -> new values for columns E, F, G, H, I
$cols = @("E", "F", "G", "H", "I")
$data = @{
    10 = @(5582, 12353, 16340, 39727, 383874)
    11 = @(127146, 130117, 171477, 288071, 2352447)
    12 = @(64121, 33888, 77612, 288071, 2050023)
    13 = @(21075, 8626, 68118, 46245, 145389)
    14 = @(12168, 10943, 0, 41818, 86177)
    15 = @(3916, 960, 2360, 5663, 14623)
    16 = @(13131, 20205, 27277, 112276, 259757)
    17 = @(17569, 252708, 439494, 896393, 2071255)
    18 = @(79368, 0, 0, 0, 0)
    19 = @(71832, 52930, 94603, 370098, 2484340)
    20 = @(415908, 522730, 897281, 2088362, 9847885)
    26 = @(180, 164, 216, 279, 572)
    27 = @(410, 407, 367, 577, 1153)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}
